# daily auto push: 2026-01-29 02:49 UTC
# Insert one new data row at row 733 (2026/01/29, 木, 10:00, rank 201),
# pushing the existing rows 733-774 down to 734-775.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(733).Insert()

# Force column A to stay plain text so "2026/01/29" isn't auto-parsed
# into a date serial (matches the rest of the date column).
$ws.Range("A733").NumberFormat = "@"
$ws.Range("A733").Value = "2026/01/29"
$ws.Range("B733").Value = "木"
$ws.Range("C733").Value = 10
$ws.Range("D733").Value = 201
